$wb = $excel.ActiveWorkbook

$prelim = $wb.Worksheets.Item("Preliminary")

# Pool A: row10 Egypt vs Tunisia, row11 Philippines vs Iran
$prelim.Range("D10").Value = 0
$prelim.Range("F10").Value = 3
$prelim.Range("H10").Value = 19
$prelim.Range("J10").Value = 25
$prelim.Range("K10").Value = 18
$prelim.Range("M10").Value = 25
$prelim.Range("N10").Value = 22
$prelim.Range("P10").Value = 25

$prelim.Range("D11").Value = 2
$prelim.Range("F11").Value = 3
$prelim.Range("H11").Value = 25
$prelim.Range("J11").Value = 21
$prelim.Range("K11").Value = 21
$prelim.Range("M11").Value = 25
$prelim.Range("N11").Value = 25
$prelim.Range("P11").Value = 17
$prelim.Range("Q11").Value = 23
$prelim.Range("S11").Value = 25
$prelim.Range("T11").Value = 20
$prelim.Range("V11").Value = 22

# Pool B: row26 Finland vs South Korea, row27 France vs Argentina
$prelim.Range("D26").Value = 3
$prelim.Range("F26").Value = 1
$prelim.Range("H26").Value = 25
$prelim.Range("J26").Value = 18
$prelim.Range("K26").Value = 25
$prelim.Range("M26").Value = 23
$prelim.Range("N26").Value = 17
$prelim.Range("P26").Value = 25
$prelim.Range("Q26").Value = 25
$prelim.Range("S26").Value = 21

$prelim.Range("D27").Value = 2
$prelim.Range("F27").Value = 3
$prelim.Range("H27").Value = 26
$prelim.Range("J27").Value = 28
$prelim.Range("K27").Value = 23
$prelim.Range("M27").Value = 25
$prelim.Range("N27").Value = 25
$prelim.Range("P27").Value = 21
$prelim.Range("Q27").Value = 25
$prelim.Range("S27").Value = 20
$prelim.Range("T27").Value = 12
$prelim.Range("V27").Value = 15

# Pool F: row50 Belgium vs Algeria, row51 Italy vs Ukraine
$prelim.Range("D50").Value = 3
$prelim.Range("F50").Value = 0
$prelim.Range("H50").Value = 25
$prelim.Range("J50").Value = 22
$prelim.Range("K50").Value = 25
$prelim.Range("M50").Value = 20
$prelim.Range("N50").Value = 25
$prelim.Range("P50").Value = 12

$prelim.Range("D51").Value = 3
$prelim.Range("F51").Value = 0
$prelim.Range("H51").Value = 25
$prelim.Range("J51").Value = 21
$prelim.Range("K51").Value = 25
$prelim.Range("M51").Value = 22
$prelim.Range("N51").Value = 25
$prelim.Range("P51").Value = 18

# Pool H: row66 Czech Republic vs China, row67 Brazil vs Serbia
$prelim.Range("D66").Value = 3
$prelim.Range("F66").Value = 0
$prelim.Range("H66").Value = 26
$prelim.Range("J66").Value = 24
$prelim.Range("K66").Value = 25
$prelim.Range("M66").Value = 19
$prelim.Range("N66").Value = 25
$prelim.Range("P66").Value = 18

$prelim.Range("D67").Value = 0
$prelim.Range("F67").Value = 3
$prelim.Range("H67").Value = 22
$prelim.Range("J67").Value = 25
$prelim.Range("K67").Value = 20
$prelim.Range("M67").Value = 25
$prelim.Range("N67").Value = 22
$prelim.Range("P67").Value = 25

# Update selection on Preliminary sheet to D6, then activate Final Round sheet
$prelim.Range("D6").Select()

$finalRound = $wb.Worksheets.Item("Final Round")
$finalRound.Activate()
$finalRound.Range("D6").Select()
